# BP-1588 category all apis
#
# Adds a "category" column to the "Курсы" (Courses) sheet, between the
# existing "learningOutcome" and "type" columns, populates its header/value
# for the single course row, and makes the "Курсы" sheet the active tab
# (it was "Адреса центров" before).

$wb = $excel.ActiveWorkbook

# "Курсы" is the first worksheet in the workbook.
$ws = $wb.Worksheets.Item(1)

# Insert a new column at E, shifting the old "type"/"options" columns
# (previously E/F) one to the right (now F/G). Excel's column Insert
# copies the formatting of the column to its left, matching column D's
# width for the new column E.
$ws.Columns("E").Insert()
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Header (row 1) and value (row 2) for the new "category" column.
$ws.Range("E1").Value = "category"
$ws.Range("E2").Value = "Профориентация"

# The new value cell uses justified vertical alignment with wrapped text
# (a new cell style), matching the other wrapped/aligned cells in row 2.
$ws.Range("E2").WrapText = $true
$ws.Range("E2").VerticalAlignment = -4130   # xlJustify

# Make "Курсы" the active sheet/tab (previously "Адреса центров" was
# active), with E1:E2 selected and E1 as the active cell.
$ws.Activate()
$ws.Range("E1:E2").Select()
